# Apply the edits described by the commit: populate a few cells on Sheet1
# with the shared string "asfasf" and leave the selection on the last
# entered cell (J8), matching the authored worksheet/workbook diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "asfasf"
$ws.Range("A8").Value = "asfasf"
$ws.Range("J8").Value = "asfasf"

$ws.Range("J8").Select()
